$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the three pricing values in column D
$ws.Range("D33").Value = 289
$ws.Range("D34").Value = 330
$ws.Range("D35").Value = 352.2
